# Applies the revisions described in the commit "Add files via upload /
# Added revised documentation." to Sort Comparison Documentation.docx
#
# Strategy: every change is performed with Find.Execute on a Range that
# has been scoped (via IndexOf on the plain-text Content) to exactly the
# anchor text we intend to touch, so it is impossible to hit the wrong
# occurrence even though some of the words involved (e.g. "groupSortTester",
# "ComparisonRunner") recur elsewhere in the document. Driving the actual
# substitution through Find/Replace (rather than a bare Range.Text
# assignment) mirrors what Word itself does and collapses the touched
# span into a single run, matching the target markup.

$d = $word.ActiveDocument

function Get-Text {
    return $d.Content.Text
}

function Find-Offset([string]$needle, [int]$searchFrom) {
    $full = Get-Text
    $idx = $full.IndexOf($needle, $searchFrom)
    if ($idx -lt 0) {
        throw "Could not find anchor text: $needle (searching from $searchFrom)"
    }
    return $idx
}

function Replace-Anchor([string]$anchor, [string]$oldInner, [string]$newInner, [int]$searchFrom) {
    # Locates $anchor (unique/disambiguated via $searchFrom), then replaces the
    # $oldInner substring within it with $newInner using Find/Replace scoped to
    # the anchor's own Range.
    $start = Find-Offset $anchor $searchFrom
    $rng = $d.Range($start, $start + $anchor.Length)
    $ok = $rng.Find.Execute($oldInner, $false, $false, $false, $false, $false, $true, 1, $false, $newInner, 1)
    if (-not $ok) {
        throw "Replace failed for anchor: $anchor"
    }
}

# 1. Title: merge "Sort Comparison " + "Documentation" runs (no visible text change).
Replace-Anchor "Sort Comparison Documentation" "Sort Comparison Documentation" "Sort Comparison Documentation" 0

# 2. "...analysis report. The code for insertion..." -> "...The Java classes for insertion..."
Replace-Anchor "The code for insertion sort" "The code for insertion sort" "The Java classes for insertion sort" 0

# 3. "A table was also given containing" -> "An Excel table was also given containing"
Replace-Anchor "A table was also given containing" "A table was also given containing" "An Excel table was also given containing" 0

# 4. "Random arrays are created by ArrayMaker" -> "Arrays are created by ArrayMaker"
Replace-Anchor "Random arrays are created by" "Random arrays are created by" "Arrays are created by" 0

# 5. "...according to user specifications and..." -> "...according to user input and..."
Replace-Anchor "according to user specifications and" "according to user specifications and" "according to user input and" 0

# 6. Heading: merge "O" + "verview of Data Management and Flow" (no visible text change).
Replace-Anchor "Overview of Data Management and Flow" "Overview of Data Management and Flow" "Overview of Data Management and Flow" 0

# 7. "ComparisonRunner.groupSortTester(). This method is called" ->
#    "ComparisonRunner's groupSortTester method. This method is called"
$apos = [char]8217
Replace-Anchor "ComparisonRunner.groupSortTester(). This method is called" `
    "ComparisonRunner.groupSortTester(). This method is called" `
    "ComparisonRunner${apos}s groupSortTester method. This method is called" 0

# 8. Remove comma: "...other layouts, but requires more setup." -> "...other layouts but requires more setup."
Replace-Anchor "other layouts, but requires more setup." "other layouts, but requires more setup." "other layouts but requires more setup." 0

# 9. Hyphenate: "uses the commented out section" -> "uses the commented-out section"
Replace-Anchor "uses the commented out section" "uses the commented out section" "uses the commented-out section" 0

# 10. Subheading: merge "group" + "SortTester" -> "groupSortTester" (no visible text change).
#     "groupSortTester" also appears as a cross-reference inside two other paragraphs, so
#     anchor on the paragraph mark + tab that follows the heading, which is unique.
$cr = [char]13
$tab = [char]9
$headingAnchor = "groupSortTester${cr}${tab}This method synthesizes"
Replace-Anchor $headingAnchor $headingAnchor $headingAnchor 0

# 11. Merge the <w:tab/> run with the following "This method " run in the groupSortTester
#     paragraph (no visible text change).
$tabAnchor1 = "${tab}This method synthesizes"
Replace-Anchor $tabAnchor1 $tabAnchor1 $tabAnchor1 0

# 12. Same tab/"This method " run merge in the getWinningSort paragraph.
$tabAnchor2 = "${tab}This method takes the data from"
Replace-Anchor $tabAnchor2 $tabAnchor2 $tabAnchor2 0

Write-Output "Done"
